# Apply updated cryptocurrency market data (prices & 1h volume %) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep these "Price" cells as text so values like "314.56" are not
# auto-converted into numbers by Excel input parsing.
$textCells = @("D4","D5","D6","D7","D9","D10","D11","D12","D15","D17","D19","D20","D22","D23","D26","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D40","D42","D43","D44","D47","D48","D50","D51")
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "42.533.18"
$ws.Range("E2").Value = "  +0.09%  "

$ws.Range("D3").Value = "2.508.83"
$ws.Range("E3").Value = "  -0.82%  "

$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "314.56"
$ws.Range("E5").Value = "  +3.11%  "

$ws.Range("D6").Value = "94.07"
$ws.Range("E6").Value = "  -3.32%  "

$ws.Range("D7").Value = "0.577"
$ws.Range("E7").Value = "  -1.76%  "

$ws.Range("E8").Value = "  -0.18%  "

$ws.Range("D9").Value = "0.526"
$ws.Range("E9").Value = "  -2.37%  "

$ws.Range("D10").Value = "35.60"
$ws.Range("E10").Value = "  -2.95%  "

$ws.Range("D11").Value = "0.0806"
$ws.Range("E11").Value = "  -0.50%  "

$ws.Range("D12").Value = "7.50"
$ws.Range("E12").Value = "  -0.32%  "

$ws.Range("E13").Value = "  -3.83%  "

$ws.Range("D14").Value = "2.894.92"
$ws.Range("E14").Value = "  -0.79%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "15.10"
$ws.Range("E15").Value = "  -0.99%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.467.15"
$ws.Range("E16").Value = "  -2.38%  "

$ws.Range("D17").Value = "0.844"
$ws.Range("E17").Value = "  -2.02%  "

$ws.Range("D18").Value = "42.651.20"
$ws.Range("E18").Value = "  +0.23%  "

$ws.Range("D19").Value = "12.89"
$ws.Range("E19").Value = "  -0.32%  "

$ws.Range("D20").Value = "6.63"
$ws.Range("E20").Value = "  +2.80%  "

$ws.Range("D21").Value = "0.0₃0955"
$ws.Range("E21").Value = "  -2.32%  "

$ws.Range("D22").Value = "69.14"
$ws.Range("E22").Value = "  -2.81%  "

$ws.Range("D23").Value = "249.73"
$ws.Range("E23").Value = "  -0.68%  "

$ws.Range("E24").Value = "  +0.65%  "

$ws.Range("E25").Value = "  -1.13%  "

$ws.Range("D26").Value = "26.62"
$ws.Range("E26").Value = "  -0.94%  "

$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("E28").Value = "  +3.94%  "

$ws.Range("D29").Value = "40.67"
$ws.Range("E29").Value = "  +6.62%  "

$ws.Range("D30").Value = "10.24"
$ws.Range("E30").Value = "  -0.84%  "

$ws.Range("D31").Value = "5.96"
$ws.Range("E31").Value = "  -0.32%  "

$ws.Range("D32").Value = "156.19"
$ws.Range("E32").Value = "  +0.71%  "

$ws.Range("B33").Value = "Celestia"
$ws.Range("C33").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D33").Value = "19.18"
$ws.Range("E33").Value = "  +3.25%  "

$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "2.10"
$ws.Range("E34").Value = "  +1.55%  "

$ws.Range("D35").Value = "3.26"
$ws.Range("E35").Value = "  -1.75%  "

$ws.Range("D36").Value = "0.0780"
$ws.Range("E36").Value = "  -1.31%  "

$ws.Range("D37").Value = "2.62"
$ws.Range("E37").Value = "  -0.34%  "

$ws.Range("D38").Value = "0.111"
$ws.Range("E38").Value = "  -4.66%  "

$ws.Range("E39").Value = "  -1.30%  "

$ws.Range("D40").Value = "23.57"
$ws.Range("E40").Value = "  -1.11%  "

$ws.Range("E41").Value = "  +12.66%  "

$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.20%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0302"
$ws.Range("E43").Value = "  +0.60%  "

$ws.Range("D44").Value = "3.75"
$ws.Range("E44").Value = "  -2.83%  "

$ws.Range("E45").Value = "  -3.23%  "

$ws.Range("D46").Value = "2.014.85"
$ws.Range("E46").Value = "  -1.40%  "

$ws.Range("D47").Value = "85.03"
$ws.Range("E47").Value = "  +0.60%  "

$ws.Range("D48").Value = "8.78"
$ws.Range("E48").Value = "  -1.77%  "

$ws.Range("D49").Value = "2.750.92"
$ws.Range("E49").Value = "  -0.97%  "

$ws.Range("D50").Value = "73.07"
$ws.Range("E50").Value = "  +1.01%  "

$ws.Range("D51").Value = "102.02"
$ws.Range("E51").Value = "  +0.24%  "
